$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.754.73"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "1.829.14"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "230.61"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D8").Value = "39.38"
$ws.Range("E8").Value = "  -2.85%  "

$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").Value = "0.0987"
$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").Value = "2.093.36"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "34.788.14"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "240.10"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("E21").Value = "  +2.29%  "

$ws.Range("D22").Value = "4.67"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "172.08"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("D27").Value = "0.123"
$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("D28").Value = "17.33"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").Value = "1.51"
$ws.Range("E29").Value = "  -7.32%  "

$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("D32").Value = "3.92"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("D33").Value = "3.92"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  +6.92%  "

$ws.Range("E36").Value = "  +11.29%  "

$ws.Range("D37").Value = "0.699"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("D38").Value = "91.32"
$ws.Range("E38").Value = "  -1.89%  "

$ws.Range("D39").Value = "1.05"
$ws.Range("E39").Value = "  +6.15%  "

$ws.Range("D40").Value = "1.338.42"
$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("E43").Value = "  -2.03%  "

$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("E46").Value = "  -1.69%  "

$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").Value = "2.008.37"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Value = "0.0671"
$ws.Range("E50").Value = "  +3.90%  "

# Row 13: Chainlink -> WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.852.17"
$ws.Range("E13").Value = "  +1.87%  "

# Row 14: WrappedEther -> Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "11.30"
$ws.Range("E14").Value = "  +1.74%  "

# Row 51: Quant -> THORChain
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "3.21"
$ws.Range("E51").Value = "  +13.02%  "
